$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '44.099.18'
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.74%  '

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '2.225.68'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.60%  '

$ws.Cells.Item(4, 5).Value = '  -1.88%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '299.03'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -2.30%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '90.69'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -3.86%  '

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.558'
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -2.18%  '

$ws.Cells.Item(8, 5).Value = '  -0.53%  '

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.494'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -5.31%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '33.34'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -4.01%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.0779'
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -2.95%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '6.97'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -3.32%  '

$ws.Cells.Item(13, 5).Value = '  -0.65%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '2.566.37'
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '2.230.46'
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +0.56%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '13.39'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.30%  '

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.777'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -6.83%  '

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '43.932.06'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.49%  '

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '12.10'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.42%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '0.0₃0906'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -4.90%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '5.97'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -5.81%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '64.05'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -2.19%  '

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '235.39'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.76%  '

$ws.Cells.Item(24, 5).Value = '  -4.76%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '1.85'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -6.67%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '39.66'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +3.57%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '2.28'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +2.84%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '9.38'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -4.18%  '

$ws.Cells.Item(30, 2).Value = 'EthereumClassic'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '19.22'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -3.82%  '

$ws.Cells.Item(31, 2).Value = 'Monero'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '151.46'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.29%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '5.45'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -8.49%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '0.0765'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -3.71%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '2.50'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -5.09%  '

$ws.Cells.Item(35, 5).Value = '  -1.62%  '

$ws.Cells.Item(36, 5).Value = '  -5.73%  '

$ws.Cells.Item(37, 5).Value = '  -7.57%  '

$ws.Cells.Item(38, 5).Value = '  -6.43%  '

$ws.Cells.Item(39, 5).Value = '  +1.17%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '3.17'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -6.54%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '3.61'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -3.82%  '

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '13.48'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -10.88%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.84%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '1.801.17'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.34%  '

$ws.Cells.Item(45, 5).Value = '  +11.31%  '

$ws.Cells.Item(46, 5).Value = '  -3.96%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '67.90'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -3.38%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '94.50'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -4.24%  '

$ws.Cells.Item(49, 2).Value = 'BitcoinSV'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '73.09'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -7.16%  '

$ws.Cells.Item(50, 2).Value = 'FraxShare'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '7.79'
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -3.77%  '

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '4.62'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -5.60%  '
